# Update cryptos list: refresh Price (D) and Volume 1h (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.199.61"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.859.36"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "'239.50"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'42.19"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.24%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "'0.0693"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "'0.0990"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "2.128.99"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "'11.49"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "1.859.65"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "'0.676"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "'4.72"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "35.209.87"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "0.0₃0796"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'240.59"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'12.23"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Value = "'168.61"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "'1.90"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.27%  "
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").Value = "'17.64"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("D32").Value = "'4.00"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("D33").Value = "'1.80"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.83%  "
$ws.Range("D34").Value = "'4.01"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "'2.04"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.66%  "
$ws.Range("D36").Value = "'0.817"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.09%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.07%  "
$ws.Range("D38").Value = "'1.09"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("D40").Value = "'90.01"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "1.345.52"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "'0.0594"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.26%  "
$ws.Range("D43").Value = "'14.94"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'12.33"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +40.28%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "'6.59"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.16%  "
$ws.Range("D49").Value = "2.048.03"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("E51").Value = "  +0.43%  "
